$p = $ppt.ActivePresentation

$p.Slides.Item(1).Shapes.Item(1).TextFrame.TextRange.Runs(1).Text = "GIGABYTE FO27Q2_2024Sun Apr  7 15:36:29 2024"
$p.Slides.Item(2).Shapes.Item(1).TextFrame.TextRange.Runs(1).Text = "ASUS XG32WCMS_CSOT_2024Sun Apr  7 15:36:29 2024"
$p.Slides.Item(3).Shapes.Item(1).TextFrame.TextRange.Runs(1).Text = "ASUS_XG32WCS-CSOT-2024Sun Apr  7 15:36:29 2024"
$p.Slides.Item(4).Shapes.Item(1).TextFrame.TextRange.Runs(1).Text = "GIGABYTE M27UA_AUO_2024Sun Apr  7 15:36:29 2024"
$p.Slides.Item(5).Shapes.Item(1).TextFrame.TextRange.Runs(1).Text = "GIGABYTE M27QA_BOE_2024Sun Apr  7 15:36:29 2024"
